$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Counter PLL Calc")

# Replace the formula-driven value in C7 with a plain numeric literal.
$ws.Range("C7").Value = 104000000

# Force a full recalculation so dependent formulas (C11, C12, C13, etc.) refresh.
$excel.CalculateFullRebuild()

# Update the selection to match the post-edit cursor position (E10).
$ws.Activate()
$ws.Range("E10").Select()
